$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.864.89"
$ws.Range("E2").Value = "  +16.00%  "
$ws.Range("D3").Value = "1.660.12"
$ws.Range("E3").Value = "  +12.98%  "
$ws.Range("D4").Value = "'0.9909"
$ws.Range("E4").Value = "  -1.89%  "
$ws.Range("D5").Value = "'305.81"
$ws.Range("E5").Value = "  +10.57%  "
$ws.Range("D6").Value = "'0.9804"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("D7").Value = "'0.3729"
$ws.Range("E7").Value = "  +5.10%  "
$ws.Range("D8").Value = "'0.3437"
$ws.Range("E8").Value = "  +12.47%  "
$ws.Range("D9").Value = "'43.89"
$ws.Range("E9").Value = "  +11.67%  "
$ws.Range("D10").Value = "'1.170"
$ws.Range("E10").Value = "  +8.45%  "
$ws.Range("D11").Value = "'0.07203"
$ws.Range("E11").Value = "  +8.83%  "
$ws.Range("D12").Value = "'0.9790"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "'20.73"
$ws.Range("E13").Value = "  +14.87%  "
$ws.Range("D14").Value = "'5.995"
$ws.Range("E14").Value = "  +9.91%  "
$ws.Range("D15").Value = "'6.743"
$ws.Range("E15").Value = "  +9.39%  "
$ws.Range("D16").Value = "1.662.12"
$ws.Range("E16").Value = "  +13.06%  "
$ws.Range("D17").Value = "'0.00001098"
$ws.Range("E17").Value = "  +8.05%  "
$ws.Range("D18").Value = "'0.9773"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "'0.06705"
$ws.Range("E19").Value = "  +12.64%  "
$ws.Range("D20").Value = "'81.19"
$ws.Range("E20").Value = "  +17.95%  "
$ws.Range("D21").Value = "'16.41"
$ws.Range("E21").Value = "  +13.76%  "
$ws.Range("D22").Value = "'6.077"
$ws.Range("E22").Value = "  +11.10%  "
$ws.Range("D23").Value = "'11.94"
$ws.Range("E23").Value = "  +6.86%  "
$ws.Range("D24").Value = "23.894.97"
$ws.Range("E24").Value = "  +16.12%  "
$ws.Range("D25").Value = "'2.353"
$ws.Range("E25").Value = "  +3.52%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.697"
$ws.Range("E26").Value = "  +29.53%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'151.84"
$ws.Range("E27").Value = "  +4.45%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.52"
$ws.Range("E28").Value = "  +14.36%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "1.842.49"
$ws.Range("E29").Value = "  +13.01%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'126.03"
$ws.Range("E30").Value = "  +10.25%  "
$ws.Range("B31").Value = "HuobiToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D31").Value = "'4.071"
$ws.Range("E31").Value = "  +3.00%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.150"
$ws.Range("E32").Value = "  +25.36%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9921"
$ws.Range("E33").Value = "  +25.35%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.702"
$ws.Range("E34").Value = "  +18.13%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.08365"
$ws.Range("E35").Value = "  +5.36%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'12.23"
$ws.Range("E36").Value = "  +19.68%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'8.928"
$ws.Range("E37").Value = "  +22.93%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06321"
$ws.Range("E38").Value = "  +11.24%  "
$ws.Range("D39").Value = "'5.278"
$ws.Range("E39").Value = "  +12.33%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.273"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.02287"
$ws.Range("E41").Value = "  +12.76%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2053"
$ws.Range("E42").Value = "  +11.32%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6029"
$ws.Range("E43").Value = "  +15.85%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'0.9785"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.819"
$ws.Range("E45").Value = "  +8.70%  "
$ws.Range("D46").Value = "'13.19"
$ws.Range("E46").Value = "  +10.07%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5907"
$ws.Range("E47").Value = "  +14.77%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'126.70"
$ws.Range("E48").Value = "  +5.57%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.991"
$ws.Range("E49").Value = "  +10.88%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07076"
$ws.Range("E50").Value = "  +10.07%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'75.44"
$ws.Range("E51").Value = "  +13.12%  "
